$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.5419114481106
$ws.Range("C2").Value = 6.770441030704588
$ws.Range("E2").Value = 13.78238661090019
$ws.Range("F2").Value = 47.67747522506782
$ws.Range("G2").Value = 3.723751959191897
$ws.Range("J2").Value = 9.810868265854641
$ws.Range("K2").Value = 14.05439295015895
$ws.Range("M2").Value = 18.09113069073169
$ws.Range("N2").Value = 23.26151785992795
$ws.Range("B3").Value = 14.35525989604425
$ws.Range("C3").Value = 6.638617225961452
$ws.Range("E3").Value = 13.75403774515839
$ws.Range("F3").Value = 47.51544599213388
$ws.Range("G3").Value = 3.726746389467194
$ws.Range("J3").Value = 9.829244988885401
$ws.Range("K3").Value = 13.93892594740123
$ws.Range("M3").Value = 18.03927139632314
$ws.Range("N3").Value = 23.30176328786925
$ws.Range("B4").Value = 14.2436441214718
$ws.Range("C4").Value = 6.558546095625697
$ws.Range("E4").Value = 13.73973460943402
$ws.Range("F4").Value = 47.42639148104627
$ws.Range("G4").Value = 3.728680171035661
$ws.Range("J4").Value = 9.841223709541584
$ws.Range("K4").Value = 13.87129541032365
$ws.Range("M4").Value = 18.01148018114591
$ws.Range("N4").Value = 23.32834009154132
$ws.Range("B5").Value = 14.19896748798776
$ws.Range("C5").Value = 6.526184402790343
$ws.Range("E5").Value = 13.73469128248644
$ws.Range("F5").Value = 47.392745433339
$ws.Range("G5").Value = 3.729492224321775
$ws.Range("J5").Value = 9.84628031793757
$ws.Range("K5").Value = 13.84458368057605
$ws.Range("M5").Value = 18.00118189059829
$ws.Range("N5").Value = 23.33963943905966
$ws.Range("B6").Value = 14.19159935542168
$ws.Range("C6").Value = 6.520828396574889
$ws.Range("E6").Value = 13.73390140532196
$ws.Range("F6").Value = 47.38731880466521
$ws.Range("G6").Value = 3.729628518508251
$ws.Range("J6").Value = 9.847130554370018
$ws.Range("K6").Value = 13.8402002024421
$ws.Range("M6").Value = 17.99953410818431
$ws.Range("N6").Value = 23.34154401520484
$ws.Range("B7").Value = 14.24303825361489
$ws.Range("C7").Value = 6.558108504303675
$ws.Range("E7").Value = 13.73966340755902
$ws.Range("F7").Value = 47.42592698614505
$ws.Range("G7").Value = 3.728691025294777
$ws.Range("J7").Value = 9.841291194985674
$ws.Range("K7").Value = 13.87093169783269
$ws.Range("M7").Value = 18.01133712695954
$ws.Range("N7").Value = 23.32849057940489
$ws.Range("B8").Value = 14.47696775793463
$ws.Range("C8").Value = 6.724837317363237
$ws.Range("E8").Value = 13.77197106208619
$ws.Range("F8").Value = 47.61945597031952
$ws.Range("G8").Value = 3.724764733775527
$ws.Range("J8").Value = 9.817060443279628
$ws.Range("K8").Value = 14.01391852734818
$ws.Range("M8").Value = 18.07241420163177
$ws.Range("N8").Value = 23.27500713051692
$ws.Range("B9").Value = 14.95688471972663
$ws.Range("C9").Value = 7.056541757220152
$ws.Range("E9").Value = 13.85969864361577
$ws.Range("F9").Value = 48.0806753016502
$ws.Range("G9").Value = 3.717816669346165
$ws.Range("J9").Value = 9.775046096471746
$ws.Range("K9").Value = 14.31897991636296
$ws.Range("M9").Value = 18.22391407265197
$ws.Range("N9").Value = 23.18493741265504
$ws.Range("B10").Value = 15.31893664676371
$ws.Range("C10").Value = 7.300332800876082
$ws.Range("E10").Value = 13.93863162792469
$ws.Range("F10").Value = 48.46772997875946
$ws.Range("G10").Value = 3.713164521487083
$ws.Range("J10").Value = 9.74751140245414
$ws.Range("K10").Value = 14.55637004374318
$ws.Range("M10").Value = 18.35394491784835
$ws.Range("N10").Value = 23.12780056890674
$ws.Range("B11").Value = 15.48495880255416
$ws.Range("C11").Value = 7.410693769377412
$ws.Range("E11").Value = 13.97758971648363
$ws.Range("F11").Value = 48.65387322365807
$ws.Range("G11").Value = 3.71114524938828
$ws.Range("J11").Value = 9.735704549732089
$ws.Range("K11").Value = 14.66684496583188
$ws.Range("M11").Value = 18.41701250304274
$ws.Range("N11").Value = 23.10377195320724
$ws.Range("B12").Value = 15.54795629935876
$ws.Range("C12").Value = 7.45236244373233
$ws.Range("E12").Value = 13.99277183959295
$ws.Range("F12").Value = 48.72577055139703
$ws.Range("G12").Value = 3.710394464906577
$ws.Range("J12").Value = 9.731336644554935
$ws.Range("K12").Value = 14.70900167435986
$ws.Range("M12").Value = 18.44144328730013
$ws.Range("N12").Value = 23.09495548313913
$ws.Range("B13").Value = 15.53438391645249
$ws.Range("C13").Value = 7.4433944926149
$ws.Range("E13").Value = 13.98948314680394
$ws.Range("F13").Value = 48.71022416276843
$ws.Range("G13").Value = 3.710555544178942
$ws.Range("J13").Value = 9.732272770228628
$ws.Range("K13").Value = 14.69990870339747
$ws.Range("M13").Value = 18.43615753032724
$ws.Range("N13").Value = 23.09684169148677
$ws.Range("B14").Value = 15.49013941856021
$ws.Range("C14").Value = 7.414124582653132
$ws.Range("E14").Value = 13.97883020513077
$ws.Range("F14").Value = 48.65976025620827
$ws.Range("G14").Value = 3.711083204379725
$ws.Range("J14").Value = 9.735343134928591
$ws.Range("K14").Value = 14.67030696756061
$ws.Range("M14").Value = 18.41901153340674
$ws.Range("N14").Value = 23.10304095030908
$ws.Range("B15").Value = 15.46305330882655
$ws.Range("C15").Value = 7.396178655877915
$ws.Range("E15").Value = 13.97236062608558
$ws.Range("F15").Value = 48.62903193532449
$ws.Range("G15").Value = 3.711408215399389
$ws.Range("J15").Value = 9.737237239842637
$ws.Range("K15").Value = 14.65221596826315
$ws.Range("M15").Value = 18.40858009122027
$ws.Range("N15").Value = 23.10687499427871
$ws.Range("B16").Value = 15.30810802847632
$ws.Range("C16").Value = 7.293105682316138
$ws.Range("E16").Value = 13.93614624830253
$ws.Range("F16").Value = 48.45576442683427
$ws.Range("G16").Value = 3.713298430943694
$ws.Range("J16").Value = 9.748297449171767
$ws.Range("K16").Value = 14.54919726289733
$ws.Range("M16").Value = 18.3499008691508
$ws.Range("N16").Value = 23.12941042894426
$ws.Range("B17").Value = 15.2133487115094
$ws.Range("C17").Value = 7.229702597790105
$ws.Range("E17").Value = 13.9147050058103
$ws.Range("F17").Value = 48.35202256647432
$ws.Range("G17").Value = 3.714482807069933
$ws.Range("J17").Value = 9.755266440965981
$ws.Range("K17").Value = 14.48661027200952
$ws.Range("M17").Value = 18.31489588669274
$ws.Range("N17").Value = 23.14373825147373
$ws.Range("B18").Value = 15.15897381817821
$ws.Range("C18").Value = 7.193186462660107
$ws.Range("E18").Value = 13.90266034313272
$ws.Range("F18").Value = 48.29330334848017
$ws.Range("G18").Value = 3.71517316502001
$ws.Range("J18").Value = 9.759342495151937
$ws.Range("K18").Value = 14.4508479222315
$ws.Range("M18").Value = 18.29513146424653
$ws.Range("N18").Value = 23.15216400285146
$ws.Range("B19").Value = 15.14058732269077
$ws.Range("C19").Value = 7.180815825639733
$ws.Range("E19").Value = 13.89863192274891
$ws.Range("F19").Value = 48.27358640657487
$ws.Range("G19").Value = 3.715408479953887
$ws.Range("J19").Value = 9.760734208560255
$ws.Range("K19").Value = 14.43878100515136
$ws.Range("M19").Value = 18.2885034670499
$ws.Range("N19").Value = 23.15504854129037
$ws.Range("B20").Value = 15.22342318913843
$ws.Range("C20").Value = 7.23645730154807
$ws.Range("E20").Value = 13.91695774398237
$ws.Range("F20").Value = 48.36296798234855
$ws.Range("G20").Value = 3.714355783241134
$ws.Range("J20").Value = 9.754517577850988
$ws.Range("K20").Value = 14.49324858842179
$ws.Range("M20").Value = 18.3185840865122
$ws.Range("N20").Value = 23.14219390535996
$ws.Range("B21").Value = 15.5031320987148
$ws.Range("C21").Value = 7.422725538826142
$ws.Range("E21").Value = 13.98194765108181
$ws.Range("F21").Value = 48.6745448153175
$ws.Range("G21").Value = 3.710927841997065
$ws.Range("J21").Value = 9.734438498591706
$ws.Range("K21").Value = 14.6789932506578
$ws.Range("M21").Value = 18.42403296432865
$ws.Range("N21").Value = 23.10121240406384
$ws.Range("B22").Value = 15.68665933343726
$ws.Range("C22").Value = 7.543728691687293
$ws.Range("E22").Value = 14.02692131121779
$ws.Range("F22").Value = 48.88637070131067
$ws.Range("G22").Value = 3.708768292454554
$ws.Range("J22").Value = 9.721916446491667
$ws.Range("K22").Value = 14.80224995220463
$ws.Range("M22").Value = 18.4961393455432
$ws.Range("N22").Value = 23.07607626158026
$ws.Range("B23").Value = 15.58866147973371
$ws.Range("C23").Value = 7.479228395923918
$ws.Range("E23").Value = 14.00269255962136
$ws.Range("F23").Value = 48.77257919316886
$ws.Range("G23").Value = 3.709913517167065
$ws.Range("J23").Value = 9.728544817596111
$ws.Range("K23").Value = 14.73630685722587
$ws.Range("M23").Value = 18.45736803762238
$ws.Range("N23").Value = 23.08934102059282
$ws.Range("B24").Value = 15.21886819069556
$ws.Range("C24").Value = 7.233403698560357
$ws.Range("E24").Value = 13.91593840138891
$ws.Range("F24").Value = 48.35801668069071
$ws.Range("G24").Value = 3.714413181289054
$ws.Range("J24").Value = 9.754855922394674
$ws.Range("K24").Value = 14.49024671842579
$ws.Range("M24").Value = 18.3169155276508
$ws.Range("N24").Value = 23.14289151708393
$ws.Range("B25").Value = 14.82514208231892
$ws.Range("C25").Value = 6.966593569697924
$ws.Range("E25").Value = 13.83339312413261
$ws.Range("F25").Value = 47.9473108469379
$ws.Range("G25").Value = 3.719616431152008
$ws.Range("J25").Value = 9.785825242152134
$ws.Range("K25").Value = 14.2339908527654
$ws.Range("M25").Value = 18.17959241346478
$ws.Range("N25").Value = 23.20771723299382
